$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.62
$ws.Range("I2").Value = 6
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.25
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("AE2").Value = 19
$ws.Range("AI2").Value = 19

# Row 3 updates
$ws.Range("G3").Value = 1.53
$ws.Range("H3").Value = 3.75
$ws.Range("J3").Value = 2.1
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 10
$ws.Range("AD3").Value = 7.5
$ws.Range("AJ3").Value = 81
$ws.Range("AZ3").Value = 201
